$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9736773371696472
$ws.Range("B1").Value = 1.837080240249634
$ws.Range("C1").Value = 3.23013162612915
$ws.Range("D1").Value = 3.907824516296387
$ws.Range("E1").Value = 0.7555598616600037
